$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 514; this shifts the existing rows 514:554
# down to 515:555, preserving all of their data untouched.
$ws.Rows("514:514").Insert()

# Populate the newly inserted row 514 with the new data point.
$ws.Range("A514").Value = 9
$ws.Range("B514").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C514").Value = "Metropolitana"
$ws.Range("D514").Value = 45013
$ws.Range("E514").Value = 13
$ws.Range("F514").Value = 100112032
$ws.Range("G514").Value = "Zapallo italiano"
$ws.Range("H514").Value = "Sin especificar"
$ws.Range("I514").Value = "Primera"
$ws.Range("J514").Value = 340
$ws.Range("K514").Value = 5000
$ws.Range("L514").Value = 6000
$ws.Range("M514").Value = 5500
$ws.Range("N514").Value = "`$/caja 50 unidades"
$ws.Range("O514").Value = "Región Metropolitana"
$ws.Range("P514").Value = 110
$ws.Range("Q514").Value = 50
$ws.Range("R514").Value = "Hortaliza"

# Ensure the date cell keeps the date number format used by the rest of
# column D (style already carried over from Insert, but set explicitly
# to be safe).
$ws.Range("D514").NumberFormat = "YYYY-MM-DD HH:MM:SS"
